$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 64 (old row 64 "totals spacer" and everything
# below it shifts down by one, to 65/66/67/68).
$ws.Range("A64:G64").Insert()

# Fill in the newly inserted row 64 with a new timesheet entry
# (2014-03-12, 19:00-21:00).
$ws.Range("A64").Value = 2014
$ws.Range("B64").Value = 3
$ws.Range("C64").Value = 12
$ws.Range("D64").Value = 0.79166666666666663
$ws.Range("E64").Value = 0.875
$ws.Range("F64").Formula = "=(E64-D64)*24*60"
$ws.Range("G64").Formula = "=F64/60"

# Row 63's end time moved 30 minutes later (0.75 -> 0.77083333333333337),
# which ripples through F63/G63 and the totals below via formulas already
# in the sheet.
$ws.Range("E63").Value = 0.77083333333333337

# Restore the view state (scroll position + active selection), which also
# shifted down by a row because of the new row above it.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$ws.Range("A65").Select()
